# Add payment 71717170 (Cash) 2025-08-20T08:14:26
#
# The previously-last row (66, phone 71717172) had its phone number stored
# as a text value; once a new row is appended after it, it settles back to
# its natural numeric type - so we rewrite A66 as a genuine number here too.
#
# Columns: A=phone, B=amount, C=method, D=timestamp, E=original_amount,
#          F=discount_applied, G=final_amount, H=birthday_discount,
#          I=points_redeemed, J=reward_discount

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 66: phone number reverts to a plain number now that it's no longer
# the newest entry.
$ws.Range("A66").Value = 71717172

# Row 67: brand-new payment row. The phone number is entered as freshly
# captured text (leading apostrophe keeps it from being auto-converted to
# a number), matching how the newest row still carries its phone as text.
$ws.Range("A67").Value = "'71717170"
$ws.Range("B67").Value = "'"
$ws.Range("C67").Value = "Cash"
$ws.Range("D67").Value = "2025-08-20T08:14:26"
$ws.Range("E67").Value = 127
$ws.Range("F67").Value = "'"
$ws.Range("G67").Value = 107.95
$ws.Range("H67").Value = 19.05
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 0
